# Generate Report for handoff
# Adds two new files (7df83b2c-c8a4-4f11-8048-7a60cc7ecb76 and
# d5e2128d-de28-409d-9b8b-139f7356d312) to the localization-status report,
# flips the two pre-existing files from "Ready for handoff" to
# "In Translation" on the Overview sheet, and fills in the per-locale
# handoff rows for the two new files on the zh-cn / de-de sheets. The
# ".localization-config" row shifts two rows down on every sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: re-colour / re-underline a just-added hyperlink cell so it
# matches the workbook's existing custom "HyperLink" look
# (font color FF6495ED, single underline) instead of the theme-based
# blue that Hyperlinks.Add applies by default. 0xED9564 is the BGR
# encoding Excel's Font.Color setter expects for RGB 6495ED.
# ---------------------------------------------------------------------
function Style-HyperlinkCell($rng) {
    $rng.Font.Color = 0xED9564
    $rng.Font.Underline = $true
}

# =======================================================================
# Sheet 1: Overview
# =======================================================================
$ws1 = $wb.Worksheets.Item("Overview")

# Drop every existing hyperlink up front -- selective replace/delete of a
# single hyperlink duplicates the <hyperlink> entry, so the whole sheet is
# reset and rebuilt from scratch instead.
$ws1.Hyperlinks.Delete()

# Row 2 / 3: existing files move from "Ready for handoff" to "In Translation"
$ws1.Range("B2").Value2 = "In Translation"
$ws1.Range("C2").Value2 = "In Translation"
$ws1.Range("B3").Value2 = "In Translation"
$ws1.Range("C3").Value2 = "In Translation"

# Row 4 (new): 7df83b2c-c8a4-4f11-8048-7a60cc7ecb76.md
$ws1.Range("A4").Value2 = "7df83b2c-c8a4-4f11-8048-7a60cc7ecb76.md"
$ws1.Range("B4").Value2 = "Ready for handoff"
$ws1.Range("C4").Value2 = "Ready for handoff"

# Row 5 (new): d5e2128d-de28-409d-9b8b-139f7356d312.md
$ws1.Range("A5").Value2 = "d5e2128d-de28-409d-9b8b-139f7356d312.md"
$ws1.Range("B5").Value2 = "Ready for handoff"
$ws1.Range("C5").Value2 = "Ready for handoff"

# Row 6 (was row 4): .localization-config, shifted down
$ws1.Range("A6").Value2 = ".localization-config"
$ws1.Range("B6").Value2 = "Not to be localized"
$ws1.Range("C6").Value2 = "Not to be localized"

$h = $ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/3ec0e6a4c639d976b5aa2977f990601fad45cc7a/e2e/655f954f-6b64-44f4-8726-2253f37e1f0a.md", "", "", "655f954f-6b64-44f4-8726-2253f37e1f0a.md")
Style-HyperlinkCell $ws1.Range("A2")
$h = $ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/3ec0e6a4c639d976b5aa2977f990601fad45cc7a/e2e/ffb1fe4c-b672-4fc0-b195-dc14ddac2614.md", "", "", "ffb1fe4c-b672-4fc0-b195-dc14ddac2614.md")
Style-HyperlinkCell $ws1.Range("A3")
$h = $ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/3ec0e6a4c639d976b5aa2977f990601fad45cc7a/e2e/7df83b2c-c8a4-4f11-8048-7a60cc7ecb76.md", "", "", "7df83b2c-c8a4-4f11-8048-7a60cc7ecb76.md")
Style-HyperlinkCell $ws1.Range("A4")
$h = $ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/3ec0e6a4c639d976b5aa2977f990601fad45cc7a/e2e/d5e2128d-de28-409d-9b8b-139f7356d312.md", "", "", "d5e2128d-de28-409d-9b8b-139f7356d312.md")
Style-HyperlinkCell $ws1.Range("A5")
$h = $ws1.Hyperlinks.Add($ws1.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/3ec0e6a4c639d976b5aa2977f990601fad45cc7a/.localization-config", "", "", ".localization-config")
Style-HyperlinkCell $ws1.Range("A6")

# =======================================================================
# Sheet 2: zh-cn
# =======================================================================
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Hyperlinks.Delete()

# Row 2 / 3: Status column flips to "In Translation"
$ws2.Range("B2").Value2 = "In Translation"
$ws2.Range("B3").Value2 = "In Translation"

# Row 4 (new): 7df83b2c-c8a4-4f11-8048-7a60cc7ecb76
$ws2.Range("A4").Value2 = "7df83b2c-c8a4-4f11-8048-7a60cc7ecb76.md"
$ws2.Range("B4").Value2 = "Ready for handoff"
$ws2.Range("C4").Value2 = "7df83b2c-c8a4-4f11-8048-7a60cc7ecb76.87b2a8f4ad157c7a7802a5bba972431e555be65f.zh-cn.xlf"
$ws2.Range("D4").Value2 = "2016-02-15 02:39:52"
$ws2.Range("G4").Value2 = "0001-01-01 00:00:00"
$ws2.Range("H4").Value2 = "Include"

# Row 5 (new): d5e2128d-de28-409d-9b8b-139f7356d312
$ws2.Range("A5").Value2 = "d5e2128d-de28-409d-9b8b-139f7356d312.md"
$ws2.Range("B5").Value2 = "Ready for handoff"
$ws2.Range("C5").Value2 = "d5e2128d-de28-409d-9b8b-139f7356d312.bbecaa647d1f8571df28a9708d69d713afe78aba.zh-cn.xlf"
$ws2.Range("D5").Value2 = "2016-02-15 02:39:52"
$ws2.Range("G5").Value2 = "0001-01-01 00:00:00"
$ws2.Range("H5").Value2 = "Include"

# Row 6 (was row 4): .localization-config, shifted down
$ws2.Range("A6").Value2 = ".localization-config"
$ws2.Range("B6").Value2 = "Not to be localized"
$ws2.Range("D6").Value2 = "0001-01-01 00:00:00"
$ws2.Range("G6").Value2 = "0001-01-01 00:00:00"
$ws2.Range("H6").Value2 = "Ignored"

# Date-ish text columns keep the document's custom datetime number format
$ws2.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("D5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("D6").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$h = $ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/3ec0e6a4c639d976b5aa2977f990601fad45cc7a/e2e/655f954f-6b64-44f4-8726-2253f37e1f0a.md", "", "", "655f954f-6b64-44f4-8726-2253f37e1f0a.md")
Style-HyperlinkCell $ws2.Range("A2")
$h = $ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a285ddc1dee3b72bc1143d88bf48e68dffbde126/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/655f954f-6b64-44f4-8726-2253f37e1f0a.5a9438de921d115e2b6e3daacd682bc0bfdc9a8d.zh-cn.xlf", "", "", "655f954f-6b64-44f4-8726-2253f37e1f0a.5a9438de921d115e2b6e3daacd682bc0bfdc9a8d.zh-cn.xlf")
Style-HyperlinkCell $ws2.Range("C2")
$h = $ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/3ec0e6a4c639d976b5aa2977f990601fad45cc7a/e2e/ffb1fe4c-b672-4fc0-b195-dc14ddac2614.md", "", "", "ffb1fe4c-b672-4fc0-b195-dc14ddac2614.md")
Style-HyperlinkCell $ws2.Range("A3")
$h = $ws2.Hyperlinks.Add($ws2.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a285ddc1dee3b72bc1143d88bf48e68dffbde126/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/ffb1fe4c-b672-4fc0-b195-dc14ddac2614.b0a7a56b7bc305b6014ca54c351de319c7c5ea02.zh-cn.xlf", "", "", "ffb1fe4c-b672-4fc0-b195-dc14ddac2614.b0a7a56b7bc305b6014ca54c351de319c7c5ea02.zh-cn.xlf")
Style-HyperlinkCell $ws2.Range("C3")
$h = $ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/3ec0e6a4c639d976b5aa2977f990601fad45cc7a/e2e/7df83b2c-c8a4-4f11-8048-7a60cc7ecb76.md", "", "", "7df83b2c-c8a4-4f11-8048-7a60cc7ecb76.md")
Style-HyperlinkCell $ws2.Range("A4")
$h = $ws2.Hyperlinks.Add($ws2.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a285ddc1dee3b72bc1143d88bf48e68dffbde126/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/7df83b2c-c8a4-4f11-8048-7a60cc7ecb76.87b2a8f4ad157c7a7802a5bba972431e555be65f.zh-cn.xlf", "", "", "7df83b2c-c8a4-4f11-8048-7a60cc7ecb76.87b2a8f4ad157c7a7802a5bba972431e555be65f.zh-cn.xlf")
Style-HyperlinkCell $ws2.Range("C4")
$h = $ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/3ec0e6a4c639d976b5aa2977f990601fad45cc7a/e2e/d5e2128d-de28-409d-9b8b-139f7356d312.md", "", "", "d5e2128d-de28-409d-9b8b-139f7356d312.md")
Style-HyperlinkCell $ws2.Range("A5")
$h = $ws2.Hyperlinks.Add($ws2.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a285ddc1dee3b72bc1143d88bf48e68dffbde126/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/d5e2128d-de28-409d-9b8b-139f7356d312.bbecaa647d1f8571df28a9708d69d713afe78aba.zh-cn.xlf", "", "", "d5e2128d-de28-409d-9b8b-139f7356d312.bbecaa647d1f8571df28a9708d69d713afe78aba.zh-cn.xlf")
Style-HyperlinkCell $ws2.Range("C5")
$h = $ws2.Hyperlinks.Add($ws2.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/3ec0e6a4c639d976b5aa2977f990601fad45cc7a/.localization-config", "", "", ".localization-config")
Style-HyperlinkCell $ws2.Range("A6")

# =======================================================================
# Sheet 3: de-de
# =======================================================================
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Hyperlinks.Delete()

# Row 2 / 3: Status column flips to "In Translation"
$ws3.Range("B2").Value2 = "In Translation"
$ws3.Range("B3").Value2 = "In Translation"

# Row 4 (new): 7df83b2c-c8a4-4f11-8048-7a60cc7ecb76
$ws3.Range("A4").Value2 = "7df83b2c-c8a4-4f11-8048-7a60cc7ecb76.md"
$ws3.Range("B4").Value2 = "Ready for handoff"
$ws3.Range("C4").Value2 = "7df83b2c-c8a4-4f11-8048-7a60cc7ecb76.87b2a8f4ad157c7a7802a5bba972431e555be65f.de-de.xlf"
$ws3.Range("D4").Value2 = "2016-02-15 02:40:06"
$ws3.Range("G4").Value2 = "0001-01-01 00:00:00"
$ws3.Range("H4").Value2 = "Include"

# Row 5 (new): d5e2128d-de28-409d-9b8b-139f7356d312
$ws3.Range("A5").Value2 = "d5e2128d-de28-409d-9b8b-139f7356d312.md"
$ws3.Range("B5").Value2 = "Ready for handoff"
$ws3.Range("C5").Value2 = "d5e2128d-de28-409d-9b8b-139f7356d312.bbecaa647d1f8571df28a9708d69d713afe78aba.de-de.xlf"
$ws3.Range("D5").Value2 = "2016-02-15 02:40:06"
$ws3.Range("G5").Value2 = "0001-01-01 00:00:00"
$ws3.Range("H5").Value2 = "Include"

# Row 6 (was row 4): .localization-config, shifted down
$ws3.Range("A6").Value2 = ".localization-config"
$ws3.Range("B6").Value2 = "Not to be localized"
$ws3.Range("D6").Value2 = "0001-01-01 00:00:00"
$ws3.Range("G6").Value2 = "0001-01-01 00:00:00"
$ws3.Range("H6").Value2 = "Ignored"

# Date-ish text columns keep the document's custom datetime number format
$ws3.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("D5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("D6").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$h = $ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/3ec0e6a4c639d976b5aa2977f990601fad45cc7a/e2e/655f954f-6b64-44f4-8726-2253f37e1f0a.md", "", "", "655f954f-6b64-44f4-8726-2253f37e1f0a.md")
Style-HyperlinkCell $ws3.Range("A2")
$h = $ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/47c4e2e2179b82e836ab47f9cd5e48c7c0a0e0cd/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/655f954f-6b64-44f4-8726-2253f37e1f0a.5a9438de921d115e2b6e3daacd682bc0bfdc9a8d.de-de.xlf", "", "", "655f954f-6b64-44f4-8726-2253f37e1f0a.5a9438de921d115e2b6e3daacd682bc0bfdc9a8d.de-de.xlf")
Style-HyperlinkCell $ws3.Range("C2")
$h = $ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/3ec0e6a4c639d976b5aa2977f990601fad45cc7a/e2e/ffb1fe4c-b672-4fc0-b195-dc14ddac2614.md", "", "", "ffb1fe4c-b672-4fc0-b195-dc14ddac2614.md")
Style-HyperlinkCell $ws3.Range("A3")
$h = $ws3.Hyperlinks.Add($ws3.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/47c4e2e2179b82e836ab47f9cd5e48c7c0a0e0cd/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/ffb1fe4c-b672-4fc0-b195-dc14ddac2614.b0a7a56b7bc305b6014ca54c351de319c7c5ea02.de-de.xlf", "", "", "ffb1fe4c-b672-4fc0-b195-dc14ddac2614.b0a7a56b7bc305b6014ca54c351de319c7c5ea02.de-de.xlf")
Style-HyperlinkCell $ws3.Range("C3")
$h = $ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/3ec0e6a4c639d976b5aa2977f990601fad45cc7a/e2e/7df83b2c-c8a4-4f11-8048-7a60cc7ecb76.md", "", "", "7df83b2c-c8a4-4f11-8048-7a60cc7ecb76.md")
Style-HyperlinkCell $ws3.Range("A4")
$h = $ws3.Hyperlinks.Add($ws3.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/47c4e2e2179b82e836ab47f9cd5e48c7c0a0e0cd/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/7df83b2c-c8a4-4f11-8048-7a60cc7ecb76.87b2a8f4ad157c7a7802a5bba972431e555be65f.de-de.xlf", "", "", "7df83b2c-c8a4-4f11-8048-7a60cc7ecb76.87b2a8f4ad157c7a7802a5bba972431e555be65f.de-de.xlf")
Style-HyperlinkCell $ws3.Range("C4")
$h = $ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/3ec0e6a4c639d976b5aa2977f990601fad45cc7a/e2e/d5e2128d-de28-409d-9b8b-139f7356d312.md", "", "", "d5e2128d-de28-409d-9b8b-139f7356d312.md")
Style-HyperlinkCell $ws3.Range("A5")
$h = $ws3.Hyperlinks.Add($ws3.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/47c4e2e2179b82e836ab47f9cd5e48c7c0a0e0cd/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/d5e2128d-de28-409d-9b8b-139f7356d312.bbecaa647d1f8571df28a9708d69d713afe78aba.de-de.xlf", "", "", "d5e2128d-de28-409d-9b8b-139f7356d312.bbecaa647d1f8571df28a9708d69d713afe78aba.de-de.xlf")
Style-HyperlinkCell $ws3.Range("C5")
$h = $ws3.Hyperlinks.Add($ws3.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/3ec0e6a4c639d976b5aa2977f990601fad45cc7a/.localization-config", "", "", ".localization-config")
Style-HyperlinkCell $ws3.Range("A6")
